$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 743152.8
$ws.Cells.Item(17, 10).Value = 1031196.56
$ws.Cells.Item(17, 12).Value = 3093589.68
$ws.Cells.Item(17, 14).Value = -3093925.68

$ws.Cells.Item(80, 8).Value = 503.27274
$ws.Cells.Item(80, 9).Value = 169.6
$ws.Cells.Item(80, 10).Value = 781.3333
$ws.Cells.Item(80, 11).Value = 508.8
$ws.Cells.Item(80, 12).Value = 2343.9999
$ws.Cells.Item(80, 13).Value = 489.2
$ws.Cells.Item(80, 14).Value = -4339.9999

$ws.Cells.Item(83, 8).Value = 503.27274
$ws.Cells.Item(83, 9).Value = 169.6
$ws.Cells.Item(83, 10).Value = 781.3333
$ws.Cells.Item(83, 11).Value = 1526.4
$ws.Cells.Item(83, 12).Value = 7031.9997
$ws.Cells.Item(83, 13).Value = 3465.6
$ws.Cells.Item(83, 14).Value = -17015.9997

$ws.Cells.Item(92, 8).Value = 6897197.5
$ws.Cells.Item(92, 9).Value = 7693027
$ws.Cells.Item(92, 11).Value = 7693027
$ws.Cells.Item(92, 13).Value = -7691779

$ws.Cells.Item(116, 8).Value = 2997.4285
$ws.Cells.Item(116, 9).Value = 2517.3333
$ws.Cells.Item(116, 11).Value = 2517.3333
$ws.Cells.Item(116, 13).Value = 924.6667000000002

$ws.Cells.Item(141, 8).Value = 1200
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 13).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7505.719
$ws.Cells.Item(32, 9).Value = 4937.863
$ws.Cells.Item(32, 10).Value = 29332.5
$ws.Cells.Item(32, 11).Value = 4937.863
$ws.Cells.Item(32, 12).Value = 29332.5
$ws.Cells.Item(32, 13).Value = -4650.863
$ws.Cells.Item(32, 14).Value = -29906.5

$ws.Cells.Item(45, 8).Value = 22344.4
$ws.Cells.Item(45, 9).Value = 22344.4
$ws.Cells.Item(45, 11).Value = 22344.4
$ws.Cells.Item(45, 13).Value = -21967.4

$ws.Cells.Item(82, 8).Value = 0
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 14).Value = 0
$ws.Cells.Item(82, 12).ClearContents()

$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 14).Value = 0
$ws.Cells.Item(85, 12).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2255.425
$ws.Cells.Item(86, 9).Value = 2100.639
$ws.Cells.Item(86, 11).Value = 2100.639
$ws.Cells.Item(86, 13).Value = -977.6390000000001

$ws.Cells.Item(89, 8).Value = 2255.425
$ws.Cells.Item(89, 9).Value = 2100.639
$ws.Cells.Item(89, 11).Value = 10503.195
$ws.Cells.Item(89, 13).Value = -4887.195

$ws.Cells.Item(107, 8).Value = 719.5
$ws.Cells.Item(107, 9).Value = 673.0909
$ws.Cells.Item(107, 11).Value = 673.0909
$ws.Cells.Item(107, 13).Value = 1246.9091

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 4875.6665
$ws.Cells.Item(86, 9).Value = 5066.2
$ws.Cells.Item(86, 11).Value = 5066.2
$ws.Cells.Item(86, 13).Value = -3943.2

$ws.Cells.Item(89, 8).Value = 4875.6665
$ws.Cells.Item(89, 9).Value = 5066.2
$ws.Cells.Item(89, 11).Value = 25331
$ws.Cells.Item(89, 13).Value = -19715

$ws.Cells.Item(93, 8).Value = 5499.1665
$ws.Cells.Item(93, 9).Value = 5499.1665
$ws.Cells.Item(93, 11).Value = 5499.1665
$ws.Cells.Item(93, 13).Value = -3627.1665

$ws.Cells.Item(107, 8).Value = 866.44446
$ws.Cells.Item(107, 10).Value = 1058.5834
$ws.Cells.Item(107, 12).Value = 1058.5834
$ws.Cells.Item(107, 14).Value = -4898.5834

$ws.Cells.Item(134, 8).Value = 9826.213
$ws.Cells.Item(134, 9).Value = 4257.9287
$ws.Cells.Item(134, 10).Value = 56599.8
$ws.Cells.Item(134, 11).Value = 12773.7861
$ws.Cells.Item(134, 12).Value = 169799.4
$ws.Cells.Item(134, 13).Value = -10238.7861
$ws.Cells.Item(134, 14).Value = -174869.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(75, 8).Value = 400
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 400
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 13).Value = 1200
$ws.Cells.Item(75, 14).Value = -3196
$ws.Cells.Item(75, 12).ClearContents()

$ws.Cells.Item(78, 8).Value = 400
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 400
$ws.Cells.Item(78, 11).Value = 0
$ws.Cells.Item(78, 13).Value = 3600
$ws.Cells.Item(78, 14).Value = -13584
$ws.Cells.Item(78, 12).ClearContents()

$ws.Cells.Item(129, 8).Value = 1702.2727
$ws.Cells.Item(129, 9).Value = 1164.75
$ws.Cells.Item(129, 10).Value = 2009.4286
$ws.Cells.Item(129, 11).Value = 3494.25
$ws.Cells.Item(129, 12).Value = 6028.2858
$ws.Cells.Item(129, 13).Value = 1505.75
$ws.Cells.Item(129, 14).Value = -16028.2858

$ws.Cells.Item(131, 8).Value = 37409.07
$ws.Cells.Item(131, 9).Value = 143871.42
$ws.Cells.Item(131, 10).Value = 1921.619
$ws.Cells.Item(131, 11).Value = 431614.26
$ws.Cells.Item(131, 12).Value = 5764.857
$ws.Cells.Item(131, 13).Value = -426574.26
$ws.Cells.Item(131, 14).Value = -15844.857

$ws.Cells.Item(132, 8).Value = 1361.7567
$ws.Cells.Item(132, 10).Value = 1682.3334
$ws.Cells.Item(132, 12).Value = 15141.0006
$ws.Cells.Item(132, 14).Value = -20201.0006

$ws.Cells.Item(133, 8).Value = 24983.166
$ws.Cells.Item(133, 9).Value = 36966.332
$ws.Cells.Item(133, 11).Value = 110898.996
$ws.Cells.Item(133, 13).Value = -105838.996

$ws.Cells.Item(136, 8).Value = 503333.16
$ws.Cells.Item(136, 9).Value = 589509.5600000001
$ws.Cells.Item(136, 10).Value = 15000
$ws.Cells.Item(136, 11).Value = 1768528.68
$ws.Cells.Item(136, 12).Value = 45000
$ws.Cells.Item(136, 13).Value = -1763428.68
$ws.Cells.Item(136, 14).Value = -55200

$ws.Cells.Item(137, 8).Value = 3889.375
$ws.Cells.Item(137, 10).Value = 4750
$ws.Cells.Item(137, 12).Value = 14250
$ws.Cells.Item(137, 14).Value = -24450

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(34, 8).Value = 37161.5
$ws.Cells.Item(34, 10).Value = 37161.5
$ws.Cells.Item(34, 12).Value = 37161.5
$ws.Cells.Item(34, 14).Value = -37697.5

$ws.Cells.Item(76, 8).Value = 37161.5
$ws.Cells.Item(76, 10).Value = 37161.5
$ws.Cells.Item(76, 12).Value = 37161.5
$ws.Cells.Item(76, 14).Value = -37791.5

$ws.Cells.Item(79, 8).Value = 37161.5
$ws.Cells.Item(79, 10).Value = 37161.5
$ws.Cells.Item(79, 12).Value = 37161.5
$ws.Cells.Item(79, 14).Value = -39345.5

$ws.Cells.Item(113, 8).Value = 2441.0588
$ws.Cells.Item(113, 9).Value = 2539.8667
$ws.Cells.Item(113, 11).Value = 2539.8667
$ws.Cells.Item(113, 13).Value = -369.8667

$ws.Cells.Item(132, 8).Value = 3311.1
$ws.Cells.Item(132, 9).Value = 2979
$ws.Cells.Item(132, 11).Value = 8937
$ws.Cells.Item(132, 13).Value = -6407

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 369.5
$ws.Cells.Item(9, 9).Value = 389
$ws.Cells.Item(9, 10).Value = 350
$ws.Cells.Item(9, 11).Value = 389
$ws.Cells.Item(9, 12).Value = 350
$ws.Cells.Item(9, 13).Value = -165
$ws.Cells.Item(9, 14).Value = -798

$ws.Cells.Item(40, 8).Value = 4439.343
$ws.Cells.Item(40, 9).Value = 3613.1738
$ws.Cells.Item(40, 11).Value = 3613.1738
$ws.Cells.Item(40, 13).Value = -3477.1738

$ws.Cells.Item(68, 8).Value = 3383.1
$ws.Cells.Item(68, 9).Value = 3380.7144
$ws.Cells.Item(68, 10).Value = 3388.6667
$ws.Cells.Item(68, 11).Value = 3380.7144
$ws.Cells.Item(68, 12).Value = 3388.6667
$ws.Cells.Item(68, 13).Value = -2631.7144
$ws.Cells.Item(68, 14).Value = -4886.6667

$ws.Cells.Item(71, 8).Value = 3383.1
$ws.Cells.Item(71, 9).Value = 3380.7144
$ws.Cells.Item(71, 10).Value = 3388.6667
$ws.Cells.Item(71, 11).Value = 16903.572
$ws.Cells.Item(71, 12).Value = 16943.3335
$ws.Cells.Item(71, 13).Value = -13159.572
$ws.Cells.Item(71, 14).Value = -24431.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 9869
$ws.Cells.Item(62, 9).Value = 4441.5
$ws.Cells.Item(62, 10).Value = 11419.714
$ws.Cells.Item(62, 11).Value = 4441.5
$ws.Cells.Item(62, 12).Value = 11419.714
$ws.Cells.Item(62, 13).Value = -3817.5
$ws.Cells.Item(62, 14).Value = -12667.714

$ws.Cells.Item(65, 8).Value = 9869
$ws.Cells.Item(65, 9).Value = 4441.5
$ws.Cells.Item(65, 10).Value = 11419.714
$ws.Cells.Item(65, 11).Value = 22207.5
$ws.Cells.Item(65, 12).Value = 57098.57
$ws.Cells.Item(65, 13).Value = -19087.5
$ws.Cells.Item(65, 14).Value = -63338.57

$ws.Cells.Item(107, 8).Value = 40334.84
$ws.Cells.Item(107, 9).Value = 333.94736
$ws.Cells.Item(107, 11).Value = 1001.84208
$ws.Cells.Item(107, 13).Value = 918.15792

$ws.Cells.Item(112, 8).Value = 19999.5
$ws.Cells.Item(112, 10).Value = 19999.5
$ws.Cells.Item(112, 12).Value = 19999.5
$ws.Cells.Item(112, 14).Value = -22953.5
